$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 141
$ws.Range("H141").Value = 997.5
$ws.Range("I141").Value = 997.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2992.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2187.5

$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 1972.4286
$ws.Range("I122").Value = 901.6667
$ws.Range("J122").Value = 2775.5
$ws.Range("K122").Value = 2705.0001
$ws.Range("L122").Value = 8326.5
$ws.Range("M122").Value = -255.0001000000002
$ws.Range("N122").Value = -13226.5
# Row 132
$ws.Range("H132").Value = 2831.1191
$ws.Range("I132").Value = 2236.543
$ws.Range("J132").Value = 5804
$ws.Range("K132").Value = 6709.629000000001
$ws.Range("L132").Value = 17412
$ws.Range("M132").Value = -4179.629000000001
$ws.Range("N132").Value = -22472
# Row 32
$ws.Range("H32").Value = 1669.7778
$ws.Range("I32").Value = 1279.5588
$ws.Range("J32").Value = 8303.5
$ws.Range("K32").Value = 1279.5588
$ws.Range("L32").Value = 8303.5
$ws.Range("M32").Value = -992.5588
# Row 45
$ws.Range("H45").Value = 90912850
$ws.Range("I45").Value = 142858530
$ws.Range("J45").Value = 7899.5
$ws.Range("K45").Value = 142858530
$ws.Range("L45").Value = 7899.5
$ws.Range("M45").Value = -142858153
# Row 51
$ws.Range("H51").Value = 39047
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 39047
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 39047
$ws.Range("N51").Value = -40559
# Row 58
$ws.Range("H58").Value = 38000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 38000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 38000
$ws.Range("M58").Value = ""
$ws.Range("N58").Value = -38860
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2862.9473
$ws.Range("I134").Value = 1728.4706
$ws.Range("J134").Value = 12506
$ws.Range("K134").Value = 5185.4118
$ws.Range("L134").Value = 37518
$ws.Range("M134").Value = -2650.4118
# Row 64
$ws.Range("H64").Value = 3662.4
$ws.Range("I64").Value = 2102
$ws.Range("J64").Value = 6003
$ws.Range("K64").Value = 2102
$ws.Range("L64").Value = 6003
$ws.Range("M64").Value = -1877
$ws.Range("N64").Value = -6453
# Row 67
$ws.Range("H67").Value = 3662.4
$ws.Range("I67").Value = 2102
$ws.Range("J67").Value = 6003
$ws.Range("K67").Value = 2102
$ws.Range("L67").Value = 6003
$ws.Range("M67").Value = -1322
$ws.Range("N67").Value = -7563
# Row 86
$ws.Range("H86").Value = 3153.25
$ws.Range("I86").Value = 2176.818
$ws.Range("J86").Value = 5301.4
$ws.Range("K86").Value = 2176.818
$ws.Range("L86").Value = 5301.4
$ws.Range("M86").Value = -1053.818
# Row 89
$ws.Range("H89").Value = 3153.25
$ws.Range("I89").Value = 2176.818
$ws.Range("J89").Value = 5301.4
$ws.Range("K89").Value = 10884.09
$ws.Range("L89").Value = 26507
$ws.Range("M89").Value = -5268.09

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 2795.639
$ws.Range("I132").Value = 3004.3635
$ws.Range("J132").Value = 2467.6428
$ws.Range("K132").Value = 9013.0905
$ws.Range("L132").Value = 7402.928400000001
$ws.Range("M132").Value = -6483.0905
$ws.Range("N132").Value = -12462.9284
# Row 134
$ws.Range("H134").Value = 1803.65
$ws.Range("I134").Value = 1164.6666
$ws.Range("J134").Value = 3720.6
$ws.Range("K134").Value = 3493.9998
$ws.Range("L134").Value = 11161.8
$ws.Range("M134").Value = -958.9998000000001
$ws.Range("N134").Value = -16231.8

$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 6414917
$ws.Range("I129").Value = 2465
$ws.Range("J129").Value = 8338652.5
$ws.Range("K129").Value = 7395
$ws.Range("L129").Value = 25015957.5
$ws.Range("M129").Value = -2395
$ws.Range("N129").Value = -25025957.5
# Row 34
$ws.Range("H34").Value = 2427.4546
$ws.Range("I34").Value = 2108.6667
$ws.Range("J34").Value = 2810
$ws.Range("K34").Value = 6326.000100000001
$ws.Range("L34").Value = 8430
$ws.Range("M34").Value = -6242.000100000001
$ws.Range("N34").Value = -8598
# Row 51
$ws.Range("H51").Value = 1311.8182
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 1311.8182
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3935.4546
$ws.Range("N51").Value = -4855.4546
# Row 57
$ws.Range("H57").Value = 5000
$ws.Range("I57").Value = 5000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -14441
# Row 58
$ws.Range("H58").Value = 500
$ws.Range("I58").Value = 500
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1500
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1372
$ws.Range("N58").Value = ""

$ws = $wb.Worksheets.Item("GSM")
# Row 124
$ws.Range("H124").Value = 70963.25
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 70963.25
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 70963.25
$ws.Range("N124").Value = -80783.25
# Row 132
$ws.Range("H132").Value = 5833.4
$ws.Range("I132").Value = 3373.2222
$ws.Range("J132").Value = 12159.571
$ws.Range("K132").Value = 10119.6666
$ws.Range("L132").Value = 36478.713
$ws.Range("M132").Value = -7589.6666

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 8887.223
$ws.Range("I122").Value = 8282.857
$ws.Range("J122").Value = 11002.5
$ws.Range("K122").Value = 24848.571
$ws.Range("L122").Value = 33007.5
$ws.Range("M122").Value = -22398.571
$ws.Range("N122").Value = -37907.5
# Row 126
$ws.Range("H126").Value = 8486.764999999999
$ws.Range("I126").Value = 4480.6
$ws.Range("J126").Value = 10156
$ws.Range("K126").Value = 13441.8
$ws.Range("L126").Value = 30468
$ws.Range("M126").Value = -10971.8
# Row 133
$ws.Range("H133").Value = 60324.9
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 60324.9
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 60324.9
$ws.Range("N133").Value = -65384.9
# Row 136
$ws.Range("H136").Value = 6618.8335
$ws.Range("I136").Value = 1662.9
$ws.Range("J136").Value = 12813.75
$ws.Range("K136").Value = 4988.700000000001
$ws.Range("L136").Value = 38441.25
$ws.Range("M136").Value = -2438.700000000001
$ws.Range("N136").Value = -43541.25
# Row 22
$ws.Range("H22").Value = 4062.2
$ws.Range("I22").Value = 1842
$ws.Range("J22").Value = 6282.4
$ws.Range("K22").Value = 1842
$ws.Range("L22").Value = 6282.4
$ws.Range("M22").Value = -1547
$ws.Range("N22").Value = -6872.4
# Row 27
$ws.Range("H27").Value = 4062.2
$ws.Range("I27").Value = 1842
$ws.Range("J27").Value = 6282.4
$ws.Range("K27").Value = 1842
$ws.Range("L27").Value = 6282.4
$ws.Range("M27").Value = -1735
$ws.Range("N27").Value = -6496.4
# Row 55
$ws.Range("H55").Value = 4526.923
$ws.Range("I55").Value = 1378.1666
$ws.Range("J55").Value = 7225.857
$ws.Range("K55").Value = 1378.1666
$ws.Range("L55").Value = 7225.857
$ws.Range("M55").Value = -1205.1666
$ws.Range("N55").Value = -7571.857
# Row 7
$ws.Range("H7").Value = 8486.764999999999
$ws.Range("I7").Value = 4480.6
$ws.Range("J7").Value = 10156
$ws.Range("K7").Value = 4480.6
$ws.Range("L7").Value = 10156
$ws.Range("M7").Value = -4368.6

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 8986.143
$ws.Range("I62").Value = 8301.333000000001
$ws.Range("J62").Value = 9499.75
$ws.Range("K62").Value = 8301.333000000001
$ws.Range("L62").Value = 9499.75
$ws.Range("M62").Value = -7677.333000000001
# Row 65
$ws.Range("H65").Value = 8986.143
$ws.Range("I65").Value = 8301.333000000001
$ws.Range("J65").Value = 9499.75
$ws.Range("K65").Value = 41506.665
$ws.Range("L65").Value = 47498.75
$ws.Range("M65").Value = -38386.665
